$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 78.916664
$ws.Range("I6").Value = 71.14286
$ws.Range("K6").Value = 213.42858
$ws.Range("M6").Value = -101.42858
$ws.Range("H96").Value = 696.8570999999999
$ws.Range("I96").Value = 319.6
$ws.Range("J96").Value = 1640
$ws.Range("K96").Value = 958.8000000000001
$ws.Range("L96").Value = 4920
$ws.Range("M96").Value = 414.1999999999999
$ws.Range("N96").Value = -7666
$ws.Range("H138").Value = 2386.2273
$ws.Range("I138").Value = 868.1539
$ws.Range("K138").Value = 2604.4617
$ws.Range("M138").Value = 2535.5383

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2997.3333
$ws.Range("I26").Value = 2997.3333
$ws.Range("K26").Value = 2997.3333
$ws.Range("M26").Value = -2667.3333
$ws.Range("H45").Value = 2491.8
$ws.Range("I45").Value = 1805.1666
$ws.Range("J45").Value = 3521.75
$ws.Range("K45").Value = 1805.1666
$ws.Range("L45").Value = 3521.75
$ws.Range("M45").Value = -1428.1666
$ws.Range("N45").Value = -4275.75
$ws.Range("H61").Value = 1484.1818
$ws.Range("I61").Value = 1347.6
$ws.Range("K61").Value = 1347.6
$ws.Range("M61").Value = -1135.6
$ws.Range("H74").Value = 3167.4546
$ws.Range("I74").Value = 3016.95
$ws.Range("J74").Value = 4672.5
$ws.Range("K74").Value = 3016.95
$ws.Range("L74").Value = 4672.5
$ws.Range("M74").Value = -2142.95
$ws.Range("N74").Value = -6420.5
$ws.Range("H77").Value = 3167.4546
$ws.Range("I77").Value = 3016.95
$ws.Range("J77").Value = 4672.5
$ws.Range("K77").Value = 15084.75
$ws.Range("L77").Value = 23362.5
$ws.Range("M77").Value = -10716.75
$ws.Range("N77").Value = -32098.5
$ws.Range("H136").Value = 1484.1818
$ws.Range("I136").Value = 1347.6
$ws.Range("K136").Value = 4042.8
$ws.Range("M136").Value = -1492.8

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5436.9443
$ws.Range("I31").Value = 4094.9333
$ws.Range("K31").Value = 4094.9333
$ws.Range("M31").Value = -3799.9333
$ws.Range("H33").Value = 7800
$ws.Range("I33").Value = 2000
$ws.Range("J33").Value = 9250
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 9250
$ws.Range("M33").Value = -1621
$ws.Range("N33").Value = -10008
$ws.Range("H34").Value = 5436.9443
$ws.Range("I34").Value = 4094.9333
$ws.Range("K34").Value = 4094.9333
$ws.Range("M34").Value = -3892.9333
$ws.Range("M44").ClearContents()
$ws.Range("H44").Value = 19998.334
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 19998.334
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 19998.334
$ws.Range("N44").Value = -20882.334
$ws.Range("H55").Value = 9891.286
$ws.Range("I55").Value = 8999
$ws.Range("J55").Value = 10248.2
$ws.Range("K55").Value = 8999
$ws.Range("L55").Value = 10248.2
$ws.Range("M55").Value = -8684
$ws.Range("N55").Value = -10878.2
$ws.Range("H58").Value = 4354
$ws.Range("J58").Value = 5095.8
$ws.Range("L58").Value = 5095.8
$ws.Range("N58").Value = -5501.8
$ws.Range("H132").Value = 1953.6041
$ws.Range("I132").Value = 1743.6578
$ws.Range("J132").Value = 2751.4
$ws.Range("K132").Value = 5230.9734
$ws.Range("L132").Value = 8254.200000000001
$ws.Range("M132").Value = -2700.9734
$ws.Range("N132").Value = -13314.2
$ws.Range("H136").Value = 4354
$ws.Range("J136").Value = 5095.8
$ws.Range("L136").Value = 15287.4
$ws.Range("N136").Value = -20387.4

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1590016.2
$ws.Range("I4").Value = 1590016.2
$ws.Range("K4").Value = 4770048.6
$ws.Range("M4").Value = -4769936.6
$ws.Range("H12").Value = 181.13333
$ws.Range("J12").Value = 191.28572
$ws.Range("L12").Value = 573.85716
$ws.Range("N12").Value = -919.85716
$ws.Range("H34").Value = 2678.4
$ws.Range("I34").Value = 149
$ws.Range("J34").Value = 3310.75
$ws.Range("K34").Value = 447
$ws.Range("L34").Value = 9932.25
$ws.Range("M34").Value = -363
$ws.Range("N34").Value = -10100.25
$ws.Range("H131").Value = 1970.3914
$ws.Range("I131").Value = 1428.0834
$ws.Range("J131").Value = 2562
$ws.Range("K131").Value = 4284.2502
$ws.Range("L131").Value = 7686
$ws.Range("M131").Value = 755.7497999999996
$ws.Range("N131").Value = -17766

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 96627.09
$ws.Range("I132").Value = 130462.25
$ws.Range("K132").Value = 391386.75
$ws.Range("M132").Value = -388856.75

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3174.1667
$ws.Range("I61").Value = 2261.25
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2261.25
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2059.25
$ws.Range("N61").Value = -5404
$ws.Range("H68").Value = 3556.6365
$ws.Range("I68").Value = 2642.7144
$ws.Range("J68").Value = 5156
$ws.Range("K68").Value = 2642.7144
$ws.Range("L68").Value = 5156
$ws.Range("M68").Value = -1893.7144
$ws.Range("N68").Value = -6654
$ws.Range("H71").Value = 3556.6365
$ws.Range("I71").Value = 2642.7144
$ws.Range("J71").Value = 5156
$ws.Range("K71").Value = 13213.572
$ws.Range("L71").Value = 25780
$ws.Range("M71").Value = -9469.572
$ws.Range("N71").Value = -33268
$ws.Range("M76").ClearContents()
$ws.Range("H76").Value = 14999.5
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 14999.5
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 14999.5
$ws.Range("N76").Value = -15675.5
$ws.Range("M79").ClearContents()
$ws.Range("H79").Value = 14999.5
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 14999.5
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 14999.5
$ws.Range("N79").Value = -17339.5
$ws.Range("N110").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H113").Value = 3174.1667
$ws.Range("I113").Value = 2261.25
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 2261.25
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -91.25
$ws.Range("N113").Value = -9340
$ws.Range("H136").Value = 4090.3635
$ws.Range("I136").Value = 2497
$ws.Range("K136").Value = 7491
$ws.Range("M136").Value = -4941

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 14011
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 14011
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 14011
$ws.Range("N20").Value = -14491
$ws.Range("M29").ClearContents()
$ws.Range("H29").Value = 4966.6665
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 4966.6665
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 4966.6665
$ws.Range("N29").Value = -5546.6665
$ws.Range("H34").Value = 36111
$ws.Range("I34").Value = 22222
$ws.Range("J34").Value = 50000
$ws.Range("K34").Value = 22222
$ws.Range("L34").Value = 50000
$ws.Range("M34").Value = -22019
$ws.Range("N34").Value = -50406
$ws.Range("H132").Value = 2277.4666
$ws.Range("I132").Value = 2089.3845
$ws.Range("K132").Value = 6268.1535
$ws.Range("M132").Value = -3738.1535
$ws.Range("H136").Value = 2288.2856
$ws.Range("I136").Value = 1636.6061
$ws.Range("J136").Value = 3521.75
$ws.Range("K136").Value = 4909.8183
$ws.Range("M136").Value = -2359.8183

Write-Host "Applied Halicarnassus_Profits updates"